# "Generate Report for Archive"
#
# The localization status report is being regenerated: the handoff that
# used to sit at "Ready for handoff" has since moved on to
# "In Translation". Update the status cells on every sheet that shows it
# (the Overview roll-up plus each per-locale detail sheet), then let the
# Status columns shrink to fit the new, shorter text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Columns.Item(3).ColumnWidth = 12.5
